$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TMF8801")

# The APPREV_MINOR (0x12) register entry that used to live in row 5 was a
# duplicate of the one already present further down the table (now row 22
# before the delete). Remove that duplicate row; everything below shifts up.
$ws.Rows.Item(5).Delete()

# The register map was re-derived for 8-bit-wide registers: every remaining
# row whose Bit Width (col D) was still 16 with Bit Index High (col E) of 15
# needs to become an 8-bit register spanning bits 7-0.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $width = $ws.Cells.Item($r, 4).Value2
    $high = $ws.Cells.Item($r, 5).Value2
    if ($width -eq 16 -and $high -eq 15) {
        $ws.Cells.Item($r, 4).Value = 8
        $ws.Cells.Item($r, 5).Value = 7
    }
}

# Restore the user's selection on the sheet to match where they left off.
$ws.Range("E13").Select()
